# Update crypto price table (rows 2-51) per scraped data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '21.903.13'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +6.77%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.575.30'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +7.01%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.49%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9884'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '287.63'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.98%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3701'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.36%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3293'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +7.79%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.146'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +8.09%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.87'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.15%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07043'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +6.16%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9971'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.11%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.28'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +11.84%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.877'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +7.41%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.549'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.01%  '

$ws.Range('B16').Value = 'Dai'
$ws.Range('C16').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9871'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.25%  '

$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001078'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.50%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.573.70'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +6.61%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06392'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +8.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '75.94'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.45%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.19'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +11.35%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.886'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +7.94%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.69'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +5.70%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '21.911.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.46%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.365'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.14%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.435'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +13.92%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '150.11'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.54%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.70'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +8.29%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.745.14'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +6.98%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '120.54'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +5.88%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.169'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +5.96%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9287'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +14.04%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.499'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.72%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08234'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.85%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.618'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.20%  '

$ws.Range('B36').Value = 'InternetComputer(DFINITY)'
$ws.Range('C36').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.229'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.64%  '

$ws.Range('B37').Value = 'FraxShare'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.671'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +13.43%  '

$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '11.89'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.77%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06184'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.22%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.237'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.60%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02196'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.67%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.2017'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.15%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9873'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.34%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5815'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.65%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.93'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.70%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.647'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.01%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5680'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +9.24%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.55'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.37%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.920'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.02%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06786'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.07%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.44'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +8.02%  '
